# Insert a new weekly data row for "Femacal de La Calera - Bruselas (repollito)".
# This shifts the existing rows 75-96 down to 76-97 and populates the new
# row 75 with the latest price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 75, pushing rows 75-96 down to 76-97.
$ws.Rows.Item(75).Insert()

# Populate the newly inserted row 75 with the new observation.
$ws.Range("A75").Value = 3
$ws.Range("B75").Value = "Femacal de La Calera"
$ws.Range("C75").Value = "Coquimbo"
$ws.Range("D75").Value = 44841
$ws.Range("E75").Value = 5
$ws.Range("F75").Value = 100112035
$ws.Range("G75").Value = "Bruselas (repollito)"
$ws.Range("H75").Value = "Sin especificar"
$ws.Range("I75").Value = "Primera"
$ws.Range("J75").Value = 38
$ws.Range("K75").Value = 15000
$ws.Range("L75").Value = 15000
$ws.Range("M75").Value = 15000
$ws.Range("N75").Value = "$/malla 15 kilos"
$ws.Range("O75").Value = "Provincia de Quillota"
$ws.Range("P75").Value = 1000
$ws.Range("Q75").Value = 15
$ws.Range("R75").Value = "Hortaliza"
